# Update the "Major Components" BOM list:
#  - Row 9  (was "Joystick")    -> "Joystick breakout board module" (component renamed)
#  - Row 10 (was "Motor (P)")   -> keeps name, but its source/part ref changes to "28BYJ-48 "
#  - Row 11 (was "Audrino board" / audrino parts) -> replaced entirely by
#           "Collision sensor limit switch module " component
#  - Row 12 (was "Driver board")-> "Driver board module" (component renamed)
#  - Row 13 unchanged content ("power supply ")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Major Components ")

# --- Row 9: Joystick -> Joystick breakout board module -------------------
$ws.Range("A9").Value = "Joystick breakout board module"

# --- Row 10: Motor (P) -- characteristic/part reference changes ----------
$ws.Range("G10").Value = "28BYJ-48 "

# --- Row 11: replace old Audrino-board row with the new Collision sensor
#             limit switch module row -------------------------------------
$ws.Range("A11").Value = "Collision sensor limit switch module "
$ws.Range("B11").Value = "5 v"
$ws.Range("C11").Value = "10mA"
$ws.Range("F11").Value = "generic"
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 3

# --- Row 12: Driver board -> Driver board module --------------------------
$ws.Range("A12").Value = "Driver board module"

# --- View/formatting touch-ups seen in the diff ---------------------------
$ws.Columns.Item(1).ColumnWidth = 27.333333333333332
$ws.Range("E22").Select()

$wb.Save()
